$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Contact": drop the unused/extra columns, promote the email column to
# "D", add a new "Full Name" column and rename the last column's header from
# "ContactTypesValue" to "ContactType".
# ---------------------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("Contact")

# Remove the big block of extra fields (old columns H:S).
$wsContact.Range("H1:S1").EntireColumn.Delete()

# Remove the empty spacer column D - shifts old E/F/G (Email/Phone/ContactTypesValue) left to D/E/F.
$wsContact.Columns.Item(4).Delete()

# Insert a brand-new column before the old "ContactTypesValue" column (now F)
# so it becomes the new "Full Name" column, pushing ContactTypesValue to G.
$wsContact.Columns.Item(6).Insert()
$wsContact.Range("F2:F3").ClearFormats()

$wsContact.Range("F1").Value = "Full Name"
$wsContact.Range("F2").Value = "Sample John"
$wsContact.Range("F3").Value = "HRSample Jing"

$wsContact.Range("G1").Value = "ContactType"

# Re-anchor the hyperlink that used to live on E2 onto the new D2 (the email
# column) and make sure no stray hyperlink formatting leaks onto the cell.
$wsContact.Range("E2").Hyperlinks.Delete()
$wsContact.Hyperlinks.Add($wsContact.Range("D2"), "mailto:johnSample@email.com")
$wsContact.Range("D2").ClearFormats()
$wsContact.Range("D2").Value2 = "johnSample@email.com"

# ---------------------------------------------------------------------------
# Sheet "ContactTypes": bold the header and select the whole column.
# ---------------------------------------------------------------------------
$wsContactTypes = $wb.Worksheets.Item("ContactTypes")
$wsContactTypes.Range("A1").Font.Bold = $true
$wsContactTypes.Range("A1:A1048576").Select()

# ---------------------------------------------------------------------------
# Sheet "Contact" becomes the active tab again, with column G selected.
# ---------------------------------------------------------------------------
$wsContact.Range("G1:G1048576").Select()
$wsContact.Activate()

Write-Host "done"
